$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.237.42"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "'1.704.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.78%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'223.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.72%  "
$ws.Range("D6").Value = "'0.5298"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'0.2650"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.52%  "
$ws.Range("D9").Value = "'0.06574"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.76%  "
$ws.Range("D10").Value = "'20.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.22%  "
$ws.Range("D11").Value = "'0.07643"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.92%  "
$ws.Range("D12").Value = "'4.578"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.75%  "
$ws.Range("D13").Value = "'1.710.87"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.08%  "
$ws.Range("D14").Value = "'1.939.93"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("E15").Value = "  -4.29%  "
$ws.Range("D16").Value = "'0.0" + [char]8325 + "8165"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.58%  "
$ws.Range("D17").Value = "'67.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").Value = "'27.220.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("D19").Value = "'216.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "'4.662"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("E22").Value = "  -4.45%  "
$ws.Range("D23").Value = "'5.966"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.39%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").Value = "'142.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("D26").Value = "'1.741"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.96%  "
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").Value = "'7.242"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("D29").Value = "'16.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.86%  "
$ws.Range("D30").Value = "'0.05369"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.35%  "
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").Value = "'3.505"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.35%  "
$ws.Range("D33").Value = "'3.412"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").Value = "'2.424"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("D37").Value = "'0.9449"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.77%  "
$ws.Range("D38").Value = "'0.5849"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("D39").Value = "'0.01630"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.36%  "
$ws.Range("D40").Value = "'5.846"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").Value = "'1.003"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").Value = "'1.037.84"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").Value = "'0.8384"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'100.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.34%  "
$ws.Range("D45").Value = "'1.846.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.69%  "
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("D47").Value = "'57.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.21%  "
$ws.Range("D48").Value = "'0.4488"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").Value = "'1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'0.06598"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +11.29%  "
$ws.Range("D51").Value = "'8.094"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.00%  "
